# Add a new RPG-lib entry: "PF and LF in same program" in row 50 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Column A: Language -> "RPG" (same value already used by every other row)
$ws.Range("A50").Value = "RPG"

# Column B: Code Ref / title
$ws.Range("B50").Value = "PF and LF in same program"

# Column C: Desc / code sample (multi-line)
$codeLines = @(
  "FACADBFSS  IF   E           K DISK                              ",
  "FFSSDBWRK  UF A E           K DISK                                    ",
  "FFSSDBWRKLFUF A E           K DISK    PREFIX('X')                     ",
  "F                                     RENAME(RFSSDBWRK:RLF)           ",
  "...",
  "Please reference  ZAUTOPAY/QFSSSRC (FSSGDBPTY )"
)
$code = [string]::Join("`n", $codeLines)
$ws.Range("C50").Value = $code

# Writing a long, multi-line wrapped string into C50 makes the host
# recompute the row's auto-fit height; restore the original explicit
# height so the row formatting stays exactly as before the edit.
$ws.Rows.Item(50).RowHeight = 26.25

# Reflect the cursor / selection ending on C50, as in the authored edit.
$ws.Range("C50").Select()
